$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.825.60'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.54%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.991.85'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '543.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.87'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.90%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.703'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +6.23%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.748'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.01%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.33'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +13.22%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.69'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.633.02'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.989.87'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.18'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.59'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.13%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.23%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.65%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.752.08'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '431.86'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '96.91'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.56'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.72%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.64%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.54'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.94%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.70'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.84'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.90%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.77'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.97%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +18.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.53'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +9.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.48'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.132'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '49.01'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +17.82%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '678.79'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '65.96'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.445'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0836'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.83%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.38'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.74%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.23%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.16%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.82%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.77'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.15%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.86'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +8.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.37'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.96%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000281'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.76%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '145.15'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.92%  '
